# Auto-generated edit script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.756.05"
$ws.Range("E2").Value = "  -1.32%  "
$ws.Range("D3").Value = "3.482.08"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.99%  "
$ws.Range("D7").Value = "3.482.43"
$ws.Range("E7").Value = "  -0.73%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -1.69%  "
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.20"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.381"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").Value = "4.073.87"
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.118"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.47%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000176"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.60%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.466.69"
$ws.Range("E17").Value = "  -1.11%  "
$ws.Range("D18").Value = "63.793.78"
$ws.Range("E18").Value = "  -1.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.97%  "
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "383.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.97%  "
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("D24").Value = "3.621.31"
$ws.Range("E24").Value = "  -0.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.64%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000114"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.58%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("E32").Value = "  -1.68%  "
$ws.Range("D33").Value = "3.488.92"
$ws.Range("E33").Value = "  -0.64%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.59%  "
$ws.Range("E36").Value = "  +1.33%  "
$ws.Range("E37").Value = "  +5.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.97"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.56%  "
$ws.Range("E39").Value = "  +0.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "159.16"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0796"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.809"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.39%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.87%  "
$ws.Range("B46").Value = "ONDO"
$ws.Range("C46").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.48%  "
$ws.Range("E48").Value = "  -0.82%  "
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("D50").Value = "2.417.91"
$ws.Range("E50").Value = "  +1.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.903"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.76%  "

Write-Output "done"
